# Fix the typo in the nested "for" loop on the "Here it is as java..." slide
# (the collision-detection code sample), and add the missing
# "ball1Collided = i;" statement that records which ball was hit.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(13)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# --- Fix 1: the inner loop compared against "i" instead of "j" ---
#   "    for ( int j = i + 1; i < soccerBalls.length; j++ {"
# becomes
#   "    for ( int j = i + 1; j < soccerBalls.length; j++ {"
$para2 = $tr.Paragraphs(2)
$oldFrag = " + 1; i < "
$newFrag = " + 1; j < "
$fragIdx = $para2.Text.IndexOf($oldFrag)
$fragStart = $para2.Start + $fragIdx
$fragRange = $tr.Characters($fragStart, $oldFrag.Length)
$fragRange.Text = $newFrag

# --- Fix 2: add a new line after the "if there is a collision, store..."
#     comment that actually stores the colliding ball's index ---
$commentPara = $tr.Paragraphs(5)
$insertText = "`r                ball1Collided = i;"
$commentPara.InsertAfter($insertText) | Out-Null

$newStart = $commentPara.Start + $commentPara.Length
$newLen = $insertText.Length
$newRange = $tr.Characters($newStart, $newLen)
$newRange.Font.Size = 10

$iIdx = $newRange.Text.LastIndexOf("i")
$iRange = $tr.Characters($newStart + $iIdx, 1)
$iRange.Font.Size = 10
